$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Dental", "dental clinic", 2),
    @("Mental Health Clinic", "behavioral health clinic", 1),
    @("Shelter", "congregate, homeless", 2),
    @("Hospice", "home care", 1),
    @("ASL", "ltcf", 1)
)

$row = 35
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

$ws.Range("A40").Select()
